# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-17 22:14:48
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists
# the users who recorded/updated a given attendance session, separated by
# ", ". For a specific set of rows the order of the names in that list was
# changed (re-ordered) during the sync. This script applies those exact
# re-orderings to column G for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @(
    @{ Row = 2;   Value = "System, system, backup@backdoor.com" },
    @{ Row = 3;   Value = "System, dnasr281@gmail.com" },
    @{ Row = 4;   Value = "System, backup@backdoor.com" },
    @{ Row = 5;   Value = "System, backup@backdoor.com" },
    @{ Row = 6;   Value = "System, dnasr281@gmail.com" },
    @{ Row = 7;   Value = "System, admin@admin.com" },
    @{ Row = 8;   Value = "System, backup@backdoor.com" },
    @{ Row = 28;  Value = "System, system, backup@backdoor.com" },
    @{ Row = 29;  Value = "System, dnasr281@gmail.com" },
    @{ Row = 30;  Value = "System, backup@backdoor.com" },
    @{ Row = 31;  Value = "System, backup@backdoor.com" },
    @{ Row = 32;  Value = "System, dnasr281@gmail.com" },
    @{ Row = 33;  Value = "System, admin@admin.com" },
    @{ Row = 34;  Value = "System, backup@backdoor.com" },
    @{ Row = 54;  Value = "System, system, backup@backdoor.com" },
    @{ Row = 55;  Value = "System, dnasr281@gmail.com" },
    @{ Row = 56;  Value = "System, backup@backdoor.com" },
    @{ Row = 57;  Value = "System, backup@backdoor.com" },
    @{ Row = 58;  Value = "System, dnasr281@gmail.com" },
    @{ Row = 59;  Value = "System, admin@admin.com" },
    @{ Row = 60;  Value = "System, backup@backdoor.com" },
    @{ Row = 80;  Value = "System, backup@backdoor.com" },
    @{ Row = 81;  Value = "System, backup@backdoor.com" },
    @{ Row = 82;  Value = "System, backup@backdoor.com" },
    @{ Row = 87;  Value = "dnasr281@gmail.com, admin@admin.com" },
    @{ Row = 106; Value = "System, backup@backdoor.com" },
    @{ Row = 107; Value = "System, backup@backdoor.com" },
    @{ Row = 108; Value = "System, backup@backdoor.com" },
    @{ Row = 113; Value = "dnasr281@gmail.com, admin@admin.com" },
    @{ Row = 132; Value = "System, backup@backdoor.com" },
    @{ Row = 133; Value = "System, backup@backdoor.com" },
    @{ Row = 134; Value = "System, backup@backdoor.com" },
    @{ Row = 139; Value = "dnasr281@gmail.com, admin@admin.com" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}

$wb.Save()
